$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "Record"
$ws.Range("B21").Value = "Balanço Geral"
$ws.Range("C21").Value = "Social"
$ws.Range("D21").Value = "2025-04-01T12:53"
$ws.Range("E21").Value = "Positivo"
$ws.Range("F21").Value = "Oportunidades no Mercado de Trabalho. Em Campos, são 359 vagas no Balcão de Empregos."
